$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet at the end of the workbook.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "getDataByConditionForRestfulGet"

# --- Copy header row (A:J) and first data row (A:J) from sheet1, preserving styles ---
$ws1.Range("A1:J2").Copy($ws2.Range("A1:J2"))

# Columns K:M on sheet1 (rspStatus/rspCode/rspMessage) shift right to L:N on sheet2.
$ws1.Range("K1:M2").Copy($ws2.Range("L1:N2"))

# New column K header + row2 value ("entityFields")
$ws2.Range("K1").Value = "entityFields"
$ws2.Range("K1").Style = $ws1.Range("K1").Style

$ws2.Range("K2").Value = "deviceType,MSDistribute"
$ws2.Range("K2").Style = $ws1.Range("K2").Style

# --- Row 3: copy formatting from sheet1 row 2, then overwrite content ---
$ws1.Range("A2:J2").Copy($ws2.Range("A3:J3"))
$ws1.Range("K2:M2").Copy($ws2.Range("L3:N3"))
$ws2.Range("K3").Style = $ws1.Range("K2").Style

$ws2.Range("A2").Value = "iot-connector-test-9"
$ws2.Range("F2").Value = "complexjsonentityone"
$ws2.Range("K2").Value = "deviceType,MSDistribute"
$ws2.Range("L2").Value = 200
$ws2.Range("M2").Value = 0
$ws2.Range("N2").Value = "Operate success."

$ws2.Range("A3").Value = "iot-connector-test-10"
$ws2.Range("F3").Value = "complexjsonentitytwo"
$ws2.Range("K3").Value = "Siid1,SoeEnabled,Siid2,Siid"
$ws2.Range("L3").Value = 200
$ws2.Range("M3").Value = 0
$ws2.Range("N3").Value = "Operate success."

# Column widths (approximate autofit results captured in the target file).
$ws2.Columns.Item(1).ColumnWidth = 17.78
$ws2.Columns.Item(2).ColumnWidth = 23.33
$ws2.Columns.Item(3).ColumnWidth = 8.33
$ws2.Columns.Item(4).ColumnWidth = 11.44
$ws2.Columns.Item(5).ColumnWidth = 5.11
$ws2.Columns.Item(6).ColumnWidth = 18.44
$ws2.Columns.Item(7).ColumnWidth = 5.22
$ws2.Columns.Item(8).ColumnWidth = 8.78
$ws2.Columns.Item(9).ColumnWidth = 7.66
$ws2.Columns.Item(10).ColumnWidth = 7.22
$ws2.Columns.Item(11).ColumnWidth = 21.89
$ws2.Columns.Item(12).ColumnWidth = 8.33
$ws2.Columns.Item(13).ColumnWidth = 7.33
$ws2.Columns.Item(14).ColumnWidth = 14

# --- Sheet1 view: unfreeze / refreeze at A2, select C6, make it the non-active tab ---
$ws1.Activate()
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("C6").Select()

# --- New sheet becomes the active tab, selection K3 ---
$ws2.Activate()
$ws2.Range("K3").Select()
